$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 52
$ws1.Range("F4").Value = 3563
$ws1.Range("F5").Value = 2207
$ws1.Range("F7").Value = 171
$ws1.Range("F8").Value = 70
$ws1.Range("F10").Value = 1312
$ws1.Range("F12").Value = 1834
$ws1.Range("F13").Value = 136

# Sheet "全部类型" (sheet4): update F column (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 52
$ws4.Range("F4").Value = 3563
$ws4.Range("F5").Value = 2207
$ws4.Range("F8").Value = 171
$ws4.Range("F9").Value = 70
$ws4.Range("F13").Value = 1312
$ws4.Range("F15").Value = 1834
$ws4.Range("F16").Value = 136
